$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "repayment_20250925_20250925 (5)"

$ws.Range("H2").Value = 247

$ws.Range("D3").Value = 7
$ws.Range("E3").Value = "1,379,742.00"
$ws.Range("G3").Value = "0.84"
$ws.Range("H3").Value = 710
$ws.Range("K3").Value = "5.58"

$ws.Range("H4").Value = 407

$ws.Range("D5").Value = 10
$ws.Range("E5").Value = "1,489,369.00"
$ws.Range("G5").Value = "1.08"
$ws.Range("H5").Value = 315
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = "4.58"
$ws.Range("L5").Value = "5.88"

$ws.Range("H6").Value = 275
$ws.Range("J6").Value = 1
$ws.Range("K6").Value = "10.80"
$ws.Range("L6").Value = "5.88"

$ws.Range("D7").Value = 8
$ws.Range("E7").Value = "757,364.00"
$ws.Range("G7").Value = "0.43"
$ws.Range("H7").Value = 334

$ws.Range("H8").Value = 247

$ws.Range("D9").Value = 8
$ws.Range("E9").Value = "5,390,208.00"
$ws.Range("G9").Value = "2.96"
$ws.Range("H9").Value = 940

$ws.Range("H10").Value = 484

$ws.Range("H11").Value = 263

$ws.Range("H12").Value = 692

$ws.Range("D13").Value = 12
$ws.Range("E13").Value = "2,374,379.00"
$ws.Range("G13").Value = "1.44"
$ws.Range("H13").Value = 2.794

$ws.Range("H14").Value = 465

$ws.Range("D15").Value = 6
$ws.Range("E15").Value = "2,170,520.00"
$ws.Range("G15").Value = "1.32"
$ws.Range("H15").Value = 308

$ws.Range("D16").Value = 2
$ws.Range("E16").Value = "1,200,000.00"
$ws.Range("G16").Value = "0.71"
$ws.Range("H16").Value = 1.104

$ws.Range("D17").Value = 5
$ws.Range("E17").Value = "1,914,209.00"
$ws.Range("G17").Value = "1.48"
$ws.Range("H17").Value = 356
$ws.Range("J17").Value = 1
$ws.Range("K17").Value = "3.49"
$ws.Range("L17").Value = "6.25"

$ws.Range("H18").Value = 640
$ws.Range("J18").Value = 1
$ws.Range("K18").Value = "5.42"
$ws.Range("L18").Value = "5.88"
